$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "23.519.53"
$ws.Cells.Item(2, 5).Value = "  +1.21%  "
$ws.Cells.Item(3, 4).Value = "1.653.70"
$ws.Cells.Item(3, 5).Value = "  +2.56%  "
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$ws.Cells.Item(5, 5).Value = "  +0.01%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "302.71"
$ws.Cells.Item(6, 5).Value = "  +0.09%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3836"
$ws.Cells.Item(7, 5).Value = "  +1.38%  "
$ws.Cells.Item(8, 2).Value = "OKB"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "51.12"
$ws.Cells.Item(8, 5).Value = "  -1.12%  "
$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.3595"
$ws.Cells.Item(9, 5).Value = "  +1.88%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.08208"
$ws.Cells.Item(10, 5).Value = "  +1.48%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "1.237"
$ws.Cells.Item(11, 5).Value = "  +2.75%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.003"
$ws.Cells.Item(12, 5).Value = "  +0.02%  "
$ws.Cells.Item(13, 5).Value = "  +1.25%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.476"
$ws.Cells.Item(14, 5).Value = "  +1.86%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.494"
$ws.Cells.Item(15, 5).Value = "  +3.11%  "
$ws.Cells.Item(16, 5).Value = "  +0.90%  "
$ws.Cells.Item(17, 4).Value = "1.650.45"
$ws.Cells.Item(17, 5).Value = "  +3.32%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "97.56"
$ws.Cells.Item(18, 5).Value = "  +4.04%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06988"
$ws.Cells.Item(19, 5).Value = "  +1.33%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.815"
$ws.Cells.Item(20, 5).Value = "  +5.41%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "17.63"
$ws.Cells.Item(21, 5).Value = "  +2.49%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.001"
$ws.Cells.Item(22, 5).Value = "  -0.03%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "12.66"
$ws.Cells.Item(23, 5).Value = "  +2.85%  "
$ws.Cells.Item(24, 4).Value = "23.539.05"
$ws.Cells.Item(24, 5).Value = "  +1.37%  "
$ws.Cells.Item(25, 5).Value = "  -0.81%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.011"
$ws.Cells.Item(26, 5).Value = "  -0.57%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "21.22"
$ws.Cells.Item(27, 5).Value = "  +1.82%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "152.81"
$ws.Cells.Item(28, 5).Value = "  +1.21%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.237"
$ws.Cells.Item(29, 5).Value = "  +0.05%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "133.81"
$ws.Cells.Item(30, 5).Value = "  +1.36%  "
$ws.Cells.Item(31, 2).Value = "Filecoin"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.209"
$ws.Cells.Item(31, 5).Value = "  +11.76%  "
$ws.Cells.Item(32, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(32, 4).Value = "1.838.05"
$ws.Cells.Item(32, 5).Value = "  +3.51%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "2.255"
$ws.Cells.Item(33, 5).Value = "  +7.60%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "12.15"
$ws.Cells.Item(34, 5).Value = "  +6.81%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.059"
$ws.Cells.Item(35, 5).Value = "  -0.94%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.02807"
$ws.Cells.Item(36, 5).Value = "  +3.81%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "6.126"
$ws.Cells.Item(37, 5).Value = "  +4.76%  "
$ws.Cells.Item(38, 5).Value = "  +1.97%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.08789"
$ws.Cells.Item(39, 5).Value = "  +1.24%  "
$ws.Cells.Item(40, 5).Value = "  +1.47%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "13.29"
$ws.Cells.Item(41, 5).Value = "  +11.53%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.7008"
$ws.Cells.Item(42, 5).Value = "  +2.05%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.337"
$ws.Cells.Item(43, 5).Value = "  +0.98%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "15.98"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.6520"
$ws.Cells.Item(45, 5).Value = "  +3.48%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.001"
$ws.Cells.Item(46, 5).Value = "  +0.11%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.306"
$ws.Cells.Item(47, 5).Value = "  +2.59%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.958"
$ws.Cells.Item(48, 5).Value = "  +0.35%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.07897"
$ws.Cells.Item(49, 5).Value = "  +0.36%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "128.34"
$ws.Cells.Item(50, 5).Value = "  +0.67%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.193"
$ws.Cells.Item(51, 5).Value = "  +2.15%  "
